{"js": "// Applies the \"Rachel Lowe CE301 Abstract\" edits:\n//  1. Append \"!\" to the end of the \"Project Title\" form-field result.\n//  2. Tweak \"Dr. \" -> \"Dr \" in front of the supervisor's name.\n//  3. Expand the HECC-IT abstract paragraph (add acronym explanation +\n//     rewrite the sentence describing OH-HECC/HECC-UP).\n//  4. Extend the \"demonstration games\" paragraph with a link + feature list.\n//  5. Replace the three placeholder paragraphs (\"explain it's been used\n//     for games\", \"give examples of the features it has\", \"etc\") with\n//     real written content.\n\nconst body = context.document.body;\n\n// --- 1. \"... - HECC-IT\" -> \"... - HECC-IT!\" -------------------------------\nconst titleHits = body.search(\" - HECC-IT\", { matchCase: true });\ntitleHits.load(\"items\");\nawait context.sync();\nif (titleHits.items.length > 0) {\n  titleHits.items[0].insertText(\" - HECC-IT!\", \"Replace\");\n}\nawait context.sync();\n\n// --- 2. \"Dr. Richard Bartle\" -> \"Dr Richard Bartle\" ------------------------\nconst supervisorHits = body.search(\"Dr. Richard Bartle\", { matchCase: true });\nsupervisorHits.load(\"items\");\nawait context.sync();\nif (supervisorHits.items.length > 0) {\n  supervisorHits.items[0].insertText(\"Dr Richard Bartle\", \"Replace\");\n}\nawait context.sync();\n\n// --- 3a. Insert the acronym expansion right after the first \"HECC-IT\" -----\nconst introHits = body.search(\"HECC-IT is a toolkit\", { matchCase: true });\nintroHits.load(\"items\");\nawait context.sync();\nif (introHits.items.length > 0) {\n  introHits.items[0].insertText(\n    \"HECC-IT (Hypertext Editing and Creation Code Integrated Toolkit) is a toolkit\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// --- 3b. Rewrite the rest of that sentence ---------------------------------\nconst tailHits = body.search(\n  \"(allowing the user to visualize their game as a network of connected passages), before using the 'HECC-UP' tool to convert the .hecc file into a playable hypertext game.\",\n  { matchCase: true }\n);\ntailHits.load(\"items\");\nawait context.sync();\nif (tailHits.items.length > 0) {\n  tailHits.items[0].insertText(\n    \"provided by the tool to assist them in the editing process, before converting their .hecc file into a playable hypertext game. Unlike most existing hypertext game authoring tools, which require authors to exclusively write raw source code or exclusively use a GUI, HECC-IT has been designed to allow authors to edit their games however they want (with or without the GUI), without having to go through a convoluted process of converting their drafts into different formats before using the other editing method.\",\n    \"Replace\"\n  );\n}\nawait context.sync();\n\n// --- 4. Extend \"This tool has been used to produce several demonstration\n//        games, \" with the itch.io link + feature rundown -----------------\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nfunction findParagraphByText(items, text) {\n  return items.find((p) => p.text === text);\n}\n\nconst demoPara = findParagraphByText(\n  paras.items,\n  \"This tool has been used to produce several demonstration games, \"\n);\nif (demoPara) {\n  demoPara.insertText(\n    \"some of which can be played here: https://11belowstudio.itch.io/the-hecc-it-demo. HECC-IT supports a variety of features which an author may opt to include within their hypertext games to enhance the player experience, from markdown formatting,  to conditional statements/guard conditions, and even the ability to disable the 'back' button. The games consist of clientside HTML/JavaScript code, allowing players to play the games with ease.\",\n    \"End\"\n  );\n}\nawait context.sync();\n\n// --- 5. Replace the three placeholder paragraphs ---------------------------\nconst paras2 = body.paragraphs;\nparas2.load(\"items/text\");\nawait context.sync();\n\nconst usedForGamesPara = findParagraphByText(\n  paras2.items,\n  \"explain it's been used for games\"\n);\nif (usedForGamesPara) {\n  usedForGamesPara.insertText(\n    \"The HECC-IT tool is a desktop application written in Java. Users may choose to create a new .hecc file upon launching the tool, or they may choose to open an existing .hecc file. Users may start editing the .hecc file using the 'OH-HECC' (Optional Help for HECC) GUI, which visualises the game as a network of connected passages, allows the user to edit existing passages/add more passages to the network, provides some error detection, edit metadata for the game, and also save their work. As the name implies, there is no obligation to use it, so, if a power user does not want to use it, they don't need to use it. It also allows users to quickly and painlessly export their game via the 'HECC-UP' (HECC Ultra Parser) utility.\",\n    \"Replace\"\n  );\n}\n\nconst featuresPara = findParagraphByText(\n  paras2.items,\n  \"give examples of the features it has\"\n);\nif (featuresPara) {\n  featuresPara.insertText(\n    \"As of the time of writing, I am still using this tool to produce a full hypertext game, called 'Backblast'; a murder mystery where you are the victim. During the authoring process, I have continued to add features to the tool as needed, for the benefit of all potential users of this tool. I intend to have this game finished by the start of the project open day, and I intend to publish both Backblast and the HECC-IT tool on my itch.io page, so anybody who wants to use HECC-IT, may use HECC-IT.\",\n    \"Replace\"\n  );\n}\n\nconst etcPara = findParagraphByText(paras2.items, \"etc\");\nif (etcPara) {\n  etcPara.insertText(\n    \"HECC-IT uses the 'showdown.js' markdown to html converter (https://github.com/showdownjs/showdown), within the terms of the MIT license, in order to provide support for markdown formatting. \",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# Applies the \"Rachel Lowe CE301 Abstract\" edits:\n#  1. Append \"!\" to the end of the \"Project Title\" form-field result.\n#  2. Tweak \"Dr. \" -> \"Dr \" in front of the supervisor's name.\n#  3. Expand the HECC-IT abstract paragraph (add acronym explanation +\n#     rewrite the sentence describing OH-HECC/HECC-UP).\n#  4. Extend the \"demonstration games\" paragraph with a link + feature list.\n#  5. Replace the three placeholder paragraphs (\"explain it's been used\n#     for games\", \"give examples of the features it has\", \"etc\") with\n#     real written content.\n\n$d = $word.ActiveDocument\n\n# --- 1. \"... - HECC-IT\" -> \"... - HECC-IT!\" --------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \" - HECC-IT\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Text = \" - HECC-IT!\"\n}\n\n# --- 2. \"Dr. Richard Bartle\" -> \"Dr Richard Bartle\" ------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Dr. Richard Bartle\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Text = \"Dr Richard Bartle\"\n}\n\n# --- 3a. Insert the acronym expansion right after the first \"HECC-IT\" ------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"HECC-IT is a toolkit\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Text = \"HECC-IT (Hypertext Editing and Creation Code Integrated Toolkit) is a toolkit\"\n}\n\n# --- 3b. Rewrite the rest of that sentence ----------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"(allowing the user to visualize their game as a network of connected passages), before using the 'HECC-UP' tool to convert the .hecc file into a playable hypertext game.\"\n$find.MatchCase = $true\nif ($find.Execute()) {\n    $rng.Text = \"provided by the tool to assist them in the editing process, before converting their .hecc file into a playable hypertext game. Unlike most existing hypertext game authoring tools, which require authors to exclusively write raw source code or exclusively use a GUI, HECC-IT has been designed to allow authors to edit their games however they want (with or without the GUI), without having to go through a convoluted process of converting their drafts into different formats before using the other editing method.\"\n}\n\n# --- 4/5. Walk the paragraphs and patch the ones that changed --------------\n$paras = $d.Paragraphs\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    $ptext = $p.Range.Text\n\n    if ($ptext -eq \"This tool has been used to produce several demonstration games, `r\") {\n        $p.Range.InsertAfter(\"some of which can be played here: https://11belowstudio.itch.io/the-hecc-it-demo. HECC-IT supports a variety of features which an author may opt to include within their hypertext games to enhance the player experience, from markdown formatting,  to conditional statements/guard conditions, and even the ability to disable the 'back' button. The games consist of clientside HTML/JavaScript code, allowing players to play the games with ease.\")\n    }\n    elseif ($ptext -eq \"explain it's been used for games`r\") {\n        $r = $p.Range\n        $r.MoveEnd(1, -1)\n        $r.Text = \"The HECC-IT tool is a desktop application written in Java. Users may choose to create a new .hecc file upon launching the tool, or they may choose to open an existing .hecc file. Users may start editing the .hecc file using the 'OH-HECC' (Optional Help for HECC) GUI, which visualises the game as a network of connected passages, allows the user to edit existing passages/add more passages to the network, provides some error detection, edit metadata for the game, and also save their work. As the name implies, there is no obligation to use it, so, if a power user does not want to use it, they don't need to use it. It also allows users to quickly and painlessly export their game via the 'HECC-UP' (HECC Ultra Parser) utility.\"\n    }\n    elseif ($ptext -eq \"give examples of the features it has`r\") {\n        $r = $p.Range\n        $r.MoveEnd(1, -1)\n        $r.Text = \"As of the time of writing, I am still using this tool to produce a full hypertext game, called 'Backblast'; a murder mystery where you are the victim. During the authoring process, I have continued to add features to the tool as needed, for the benefit of all potential users of this tool. I intend to have this game finished by the start of the project open day, and I intend to publish both Backblast and the HECC-IT tool on my itch.io page, so anybody who wants to use HECC-IT, may use HECC-IT.\"\n    }\n    elseif ($ptext -eq \"etc`r\") {\n        $r = $p.Range\n        $r.MoveEnd(1, -1)\n        $r.Text = \"HECC-IT uses the 'showdown.js' markdown to html converter (https://github.com/showdownjs/showdown), within the terms of the MIT license, in order to provide support for markdown formatting. \"\n    }\n}\n"}
